$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Update the "Förändrad" (C column) date for all existing data rows (2..431)
#    from 45190 (2023-09-21) to 45192 (2023-09-23).
$ws.Range("C2:C431").Value = 45192

# 2) Row 431 gains an explicit row height (ht="15" customHeight="1").
$ws.Range("A431").EntireRow.RowHeight = 15

# 3) Append a brand-new record in row 432.
$ws.Range("A432").Value = "A 44479-2023"

$ws.Range("B432").Value = 45189
$ws.Range("B432").NumberFormat = $ws.Range("B431").NumberFormat

$ws.Range("C432").Value = 45192
$ws.Range("C432").NumberFormat = $ws.Range("C431").NumberFormat

$ws.Range("D432").Value = "SÖDERMANLANDS LÄN"
$ws.Range("E432").Value = "STRÄNGNÄS"

$ws.Range("G432").Value = 4
$ws.Range("H432").Value = 0
$ws.Range("I432").Value = 0
$ws.Range("J432").Value = 0
$ws.Range("K432").Value = 0
$ws.Range("L432").Value = 0
$ws.Range("M432").Value = 0
$ws.Range("N432").Value = 0
$ws.Range("O432").Value = 0
$ws.Range("P432").Value = 0
$ws.Range("Q432").Value = 0

$ws.Range("R432").WrapText = $true
